$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> corrected notarjeta (card number), stored as zero-padded text
$cardData = @(
    "2|020001",
    "3|020002",
    "4|020003",
    "5|020004",
    "6|020005",
    "7|020006",
    "8|020007",
    "9|020008",
    "10|020009",
    "11|020010",
    "12|020011",
    "13|020012",
    "14|020013",
    "15|020014",
    "16|020015",
    "17|020016",
    "18|020017",
    "19|020018",
    "20|020019",
    "21|020020",
    "22|020021",
    "23|020022",
    "24|020023",
    "25|020024",
    "26|020025",
    "27|020026",
    "28|020027",
    "29|020028",
    "30|020029",
    "31|020030",
    "32|020031",
    "33|020032",
    "34|020033",
    "35|020034",
    "36|020035",
    "37|020036",
    "38|050001",
    "39|050002",
    "40|050003",
    "41|050004",
    "42|050005",
    "43|050006",
    "44|050007",
    "45|050008",
    "46|050009",
    "47|050010",
    "48|050011",
    "49|050012",
    "50|050013",
    "51|050014",
    "52|050015",
    "53|050016",
    "54|050017",
    "55|050018",
    "56|050019",
    "57|050020",
    "58|070001",
    "59|070002",
    "60|070003",
    "61|070004",
    "62|070005",
    "63|070006",
    "64|070007",
    "65|070008",
    "66|070009",
    "67|070010",
    "68|070011",
    "69|070012",
    "70|070013",
    "71|070014",
    "72|070015",
    "73|070016",
    "74|070017",
    "75|070018",
    "76|070019",
    "77|070020",
    "78|070021",
    "79|070022",
    "80|070023",
    "81|070024",
    "82|070025",
    "83|070026",
    "84|070027",
    "85|070028",
    "86|070029",
    "87|070030",
    "88|070031",
    "89|070032",
    "90|070033",
    "91|070034",
    "92|070035",
    "93|070036",
    "94|070037",
    "95|070038",
    "96|070039",
    "97|070040",
    "98|070041",
    "99|070042",
    "100|070043",
    "101|070044",
    "102|070045",
    "103|070046",
    "104|070047",
    "105|070048",
    "106|070049",
    "107|070050",
    "108|070051",
    "109|060001",
    "110|060002",
    "111|060003",
    "112|060004",
    "113|060005",
    "114|060006",
    "115|060007",
    "116|060008",
    "117|060009",
    "118|060010",
    "119|060011",
    "120|060012",
    "121|060013",
    "122|060014",
    "123|060015",
    "124|060016",
    "125|060017",
    "126|060018",
    "127|060019",
    "128|060020",
    "129|060021",
    "130|060022",
    "131|060023",
    "132|060024",
    "133|060025",
    "134|060026",
    "135|060027",
    "136|060028",
    "137|060029",
    "138|060030",
    "139|060031",
    "140|060032",
    "141|060033",
    "142|060034",
    "143|060035",
    "144|060036",
    "145|060037",
    "146|060038",
    "147|060039",
    "148|060040",
    "149|060041",
    "150|060042",
    "151|060043",
    "152|060044",
    "153|060045",
    "154|060046",
    "155|060047",
    "156|060048",
    "157|060049",
    "158|060050",
    "159|060051",
    "160|060052",
    "161|060053",
    "162|060054",
    "163|060055",
    "164|060056",
    "165|060057",
    "166|060058",
    "167|060059",
    "168|060060",
    "169|060061",
    "170|060062",
    "171|060063",
    "172|060064",
    "173|060065",
    "174|060066",
    "175|060067",
    "176|060068",
    "177|060069",
    "178|060070",
    "179|060071",
    "180|060072",
    "181|060073",
    "182|060074",
    "183|060075",
    "184|060076",
    "185|060077",
    "186|060078",
    "187|060079",
    "188|060080",
    "189|060081",
    "190|060082",
    "191|060083",
    "192|060084",
    "193|060085",
    "194|060086",
    "195|060087",
    "196|060088",
    "197|060089",
    "198|060090",
    "199|060091",
    "200|060092",
    "201|060093",
    "202|060094",
    "203|060095",
    "204|060096",
    "205|060097",
    "206|080001",
    "207|080002",
    "208|080003",
    "209|080004",
    "210|080005",
    "211|080006",
    "212|080007",
    "213|080008",
    "214|080009",
    "215|080010",
    "216|080011",
    "217|080012",
    "218|080013",
    "219|080014",
    "220|080015",
    "221|080016",
    "222|080017",
    "223|080018",
    "224|080019",
    "225|080020",
    "226|080021",
    "227|080022",
    "228|080023",
    "229|080024",
    "230|080025",
    "231|080026",
    "232|090124",
    "233|090125",
    "234|090126",
    "235|090127",
    "236|090128",
    "237|090129",
    "238|090130",
    "239|090131",
    "240|090132",
    "241|090133",
    "242|090134",
    "243|090135",
    "244|090136",
    "245|090137",
    "246|090138",
    "247|090139",
    "248|090140",
    "249|090141",
    "250|090142",
    "251|090143",
    "252|090144",
    "253|090145",
    "254|090146",
    "255|090147",
    "256|010001",
    "257|010002",
    "258|010003",
    "259|010004",
    "260|090152",
    "261|090153"
)

foreach ($entry in $cardData) {
    $parts = $entry.Split("|")
    $row = [int]$parts[0]
    $val = $parts[1]
    $cell = $ws.Cells.Item($row, 9)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Department id corrections (column J) for rows that were miscategorized
$ws.Cells.Item(256, 10).Value = 1
$ws.Cells.Item(257, 10).Value = 1
$ws.Cells.Item(258, 10).Value = 1
$ws.Cells.Item(259, 10).Value = 1
$ws.Cells.Item(260, 10).Value = 23
$ws.Cells.Item(261, 10).Value = 23

# Gender correction (column E) for row 261
$ws.Cells.Item(261, 5).Value = "M"

# New employee row 262
$ws.Cells.Item(262, 1).Value = 261
$ws.Cells.Item(262, 2).Value = 275
$ws.Cells.Item(262, 3).Value = "TRUE"
$ws.Cells.Item(262, 4).Value = "TRUE"
$ws.Cells.Item(262, 5).Value = "M"
$ws.Cells.Item(262, 6).Value = "Itzel"
$ws.Cells.Item(262, 7).Value = "Rodriguez"
$ws.Cells.Item(262, 8).Value = "Rosas"
$cellI262 = $ws.Cells.Item(262, 9)
$cellI262.NumberFormat = "@"
$cellI262.Value = "090154"
$ws.Cells.Item(262, 10).Value = 23
$ws.Cells.Item(262, 11).Value = "ROLE_CH"
$ws.Cells.Item(262, 12).Value = 3

# Update sheet view to reflect scroll position / active selection
$ws.Application.ActiveWindow.ScrollRow = 219
$sheetView = $ws.Application.ActiveWindow
$ws.Range("H235").Select()
